$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.024885835776341
$ws.Cells.Item(2, 4).Value = 1.027530404479445
$ws.Cells.Item(2, 5).Value = 1.025282178188
$ws.Cells.Item(2, 6).Value = 1.023375831561647
$ws.Cells.Item(2, 9).Value = 1.029398573948809
$ws.Cells.Item(2, 10).Value = 1.030057944561719
$ws.Cells.Item(2, 11).Value = 1.030349864236994
$ws.Cells.Item(2, 12).Value = 1.028108205776566
$ws.Cells.Item(2, 13).Value = 1.026207452831525
$ws.Cells.Item(2, 14).Value = 1.013949506833425

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.026247476100909
$ws.Cells.Item(3, 4).Value = 1.028804876444447
$ws.Cells.Item(3, 5).Value = 1.026450230152765
$ws.Cells.Item(3, 6).Value = 1.025380114816553
$ws.Cells.Item(3, 9).Value = 1.029617767452617
$ws.Cells.Item(3, 10).Value = 1.031056554653581
$ws.Cells.Item(3, 11).Value = 1.031430819698768
$ws.Cells.Item(3, 12).Value = 1.029082548166909
$ws.Cells.Item(3, 13).Value = 1.028015340439793
$ws.Cells.Item(3, 14).Value = 1.014289224448835

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.027126755963512
$ws.Cells.Item(4, 4).Value = 1.029628087600074
$ws.Cells.Item(4, 5).Value = 1.027204790601749
$ws.Cells.Item(4, 6).Value = 1.026674795879949
$ws.Cells.Item(4, 9).Value = 1.029757460501956
$ws.Cells.Item(4, 10).Value = 1.031700534301232
$ws.Cells.Item(4, 11).Value = 1.032128273735062
$ws.Cells.Item(4, 12).Value = 1.029711202872723
$ws.Cells.Item(4, 13).Value = 1.02918257408318
$ws.Cells.Item(4, 14).Value = 1.014508073131047

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.02749598362172
$ws.Cells.Item(5, 4).Value = 1.029973823085706
$ws.Cells.Item(5, 5).Value = 1.027521714744117
$ws.Cells.Item(5, 6).Value = 1.027218564016546
$ws.Cells.Item(5, 9).Value = 1.029815676697336
$ws.Cells.Item(5, 10).Value = 1.031970745254772
$ws.Cells.Item(5, 11).Value = 1.032421011380872
$ws.Cells.Item(5, 12).Value = 1.029975061175865
$ws.Cells.Item(5, 13).Value = 1.029672673969944
$ws.Cells.Item(5, 14).Value = 1.014599846381376

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.027557954065214
$ws.Cells.Item(6, 4).Value = 1.030031853697334
$ws.Cells.Item(6, 5).Value = 1.027574910676143
$ws.Cells.Item(6, 6).Value = 1.027309835316511
$ws.Cells.Item(6, 9).Value = 1.029825421544901
$ws.Cells.Item(6, 10).Value = 1.032016084646084
$ws.Cells.Item(6, 11).Value = 1.032470135810565
$ws.Cells.Item(6, 12).Value = 1.030019339189465
$ws.Cells.Item(6, 13).Value = 1.029754928846783
$ws.Cells.Item(6, 14).Value = 1.014615242035524

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.027131691248875
$ws.Cells.Item(7, 4).Value = 1.029632708672384
$ws.Cells.Item(7, 5).Value = 1.027209026502423
$ws.Cells.Item(7, 6).Value = 1.026682063742553
$ws.Cells.Item(7, 9).Value = 1.029758240394714
$ws.Cells.Item(7, 10).Value = 1.031704146901986
$ws.Cells.Item(7, 11).Value = 1.032132187157311
$ws.Cells.Item(7, 12).Value = 1.029714730236213
$ws.Cells.Item(7, 13).Value = 1.029189125181162
$ws.Cells.Item(7, 14).Value = 1.014509300314133

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.025346384191264
$ws.Cells.Item(8, 4).Value = 1.027961423901445
$ws.Cells.Item(8, 5).Value = 1.025677188418627
$ws.Cells.Item(8, 6).Value = 1.024053659101676
$ws.Cells.Item(8, 9).Value = 1.029473095395312
$ws.Cells.Item(8, 10).Value = 1.030395885965844
$ws.Cells.Item(8, 11).Value = 1.030715595147467
$ws.Cells.Item(8, 12).Value = 1.028437867245299
$ws.Cells.Item(8, 13).Value = 1.026818982386159
$ws.Cells.Item(8, 14).Value = 1.014064518203286

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.022186377622743
$ws.Cells.Item(9, 4).Value = 1.025004954538403
$ws.Cells.Item(9, 5).Value = 1.022968091716969
$ws.Cells.Item(9, 6).Value = 1.019404287668784
$ws.Cells.Item(9, 9).Value = 1.02895417657644
$ws.Cells.Item(9, 10).Value = 1.028073550034312
$ws.Cells.Item(9, 11).Value = 1.028203821736755
$ws.Cells.Item(9, 12).Value = 1.026173775648194
$ws.Cells.Item(9, 13).Value = 1.022621970502636
$ws.Cells.Item(9, 14).Value = 1.013273234536917

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020069770884296
$ws.Cells.Item(10, 4).Value = 1.023025866412619
$ws.Cells.Item(10, 5).Value = 1.021155084346141
$ws.Cells.Item(10, 6).Value = 1.016291693723902
$ws.Cells.Item(10, 9).Value = 1.028597068440183
$ws.Cells.Item(10, 10).Value = 1.026513525482202
$ws.Cells.Item(10, 11).Value = 1.026518469631543
$ws.Cells.Item(10, 12).Value = 1.024654574244393
$ws.Cells.Item(10, 13).Value = 1.01980921788964
$ws.Cells.Item(10, 14).Value = 1.012740540199882

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.019150792962312
$ws.Cells.Item(11, 4).Value = 1.022166884170708
$ws.Cells.Item(11, 5).Value = 1.020368305962069
$ws.Cells.Item(11, 6).Value = 1.014940584219731
$ws.Cells.Item(11, 9).Value = 1.028439766893572
$ws.Cells.Item(11, 10).Value = 1.025835137190026
$ws.Cells.Item(11, 11).Value = 1.025786038510283
$ws.Cells.Item(11, 12).Value = 1.023994342890702
$ws.Cells.Item(11, 13).Value = 1.018587555461154
$ws.Cells.Item(11, 14).Value = 1.012508624885663

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.018809062727494
$ws.Cells.Item(12, 4).Value = 1.021847507801592
$ws.Cells.Item(12, 5).Value = 1.02007579394671
$ws.Cells.Item(12, 6).Value = 1.014438200473954
$ws.Cells.Item(12, 9).Value = 1.028380934764889
$ws.Cells.Item(12, 10).Value = 1.025582712652582
$ws.Cells.Item(12, 11).Value = 1.02551357355377
$ws.Cells.Item(12, 12).Value = 1.02374873542729
$ws.Cells.Item(12, 13).Value = 1.018133198090156
$ws.Cells.Item(12, 14).Value = 1.012422290345868

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.018882382411181
$ws.Cells.Item(13, 4).Value = 1.021916029370796
$ws.Cells.Item(13, 5).Value = 1.020138550983851
$ws.Cells.Item(13, 6).Value = 1.014545987421631
$ws.Cells.Item(13, 9).Value = 1.028393572739941
$ws.Cells.Item(13, 10).Value = 1.025636878640562
$ws.Cells.Item(13, 11).Value = 1.025572036768589
$ws.Cells.Item(13, 12).Value = 1.023801435826229
$ws.Cells.Item(13, 13).Value = 1.018230685720848
$ws.Cells.Item(13, 14).Value = 1.012440818071933

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.019122553259469
$ws.Cells.Item(14, 4).Value = 1.022140490834297
$ws.Cells.Item(14, 5).Value = 1.020344132324834
$ws.Cells.Item(14, 6).Value = 1.014899067779905
$ws.Cells.Item(14, 9).Value = 1.028434912050873
$ws.Cells.Item(14, 10).Value = 1.025814280735166
$ws.Cells.Item(14, 11).Value = 1.025763524826518
$ws.Cells.Item(14, 12).Value = 1.023974048448995
$ws.Cells.Item(14, 13).Value = 1.018550010027916
$ws.Cells.Item(14, 14).Value = 1.012501492352649

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.019270479689446
$ws.Cells.Item(15, 4).Value = 1.022278747349664
$ws.Cells.Item(15, 5).Value = 1.020470762119687
$ws.Cells.Item(15, 6).Value = 1.015116542522998
$ws.Cells.Item(15, 9).Value = 1.028460329059602
$ws.Cells.Item(15, 10).Value = 1.025923525379884
$ws.Cells.Item(15, 11).Value = 1.025881452727862
$ws.Cells.Item(15, 12).Value = 1.024080351804786
$ws.Cells.Item(15, 13).Value = 1.018746679188765
$ws.Cells.Item(15, 14).Value = 1.012538850418804

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020130707518982
$ws.Cells.Item(16, 4).Value = 1.023082830907606
$ws.Cells.Item(16, 5).Value = 1.021207263107146
$ws.Cells.Item(16, 6).Value = 1.016381290398338
$ws.Cells.Item(16, 9).Value = 1.028607451574067
$ws.Cells.Item(16, 10).Value = 1.026558486438954
$ws.Cells.Item(16, 11).Value = 1.026567021945572
$ws.Cells.Item(16, 12).Value = 1.024698340356552
$ws.Cells.Item(16, 13).Value = 1.019890215620524
$ws.Cells.Item(16, 14).Value = 1.012755905015737

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.02066963698858
$ws.Cells.Item(17, 4).Value = 1.023586663747997
$ws.Cells.Item(17, 5).Value = 1.021668781537781
$ws.Cells.Item(17, 6).Value = 1.017173725939695
$ws.Cells.Item(17, 9).Value = 1.028699021044396
$ws.Cells.Item(17, 10).Value = 1.026956002699088
$ws.Cells.Item(17, 11).Value = 1.02699634335984
$ws.Cells.Item(17, 12).Value = 1.025085339227503
$ws.Cells.Item(17, 13).Value = 1.020606516426853
$ws.Cells.Item(17, 14).Value = 1.012891719979417

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.020983747435139
$ws.Cells.Item(18, 4).Value = 1.023880346178307
$ws.Cells.Item(18, 5).Value = 1.021937810575859
$ws.Cells.Item(18, 6).Value = 1.01763561971441
$ws.Cells.Item(18, 9).Value = 1.028752174354963
$ws.Cells.Item(18, 10).Value = 1.027187589105093
$ws.Cells.Item(18, 11).Value = 1.027246502886193
$ws.Cells.Item(18, 12).Value = 1.02531083731244
$ws.Cells.Item(18, 13).Value = 1.021023964700738
$ws.Cells.Item(18, 14).Value = 1.012970817602656

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.021090810892908
$ws.Cells.Item(19, 4).Value = 1.023980451507513
$ws.Cells.Item(19, 5).Value = 1.022029514432783
$ws.Cells.Item(19, 6).Value = 1.017793059810588
$ws.Cells.Item(19, 9).Value = 1.028770254624369
$ws.Cells.Item(19, 10).Value = 1.027266507127891
$ws.Cells.Item(19, 11).Value = 1.027331757554335
$ws.Cells.Item(19, 12).Value = 1.025387687208491
$ws.Cells.Item(19, 13).Value = 1.021166243626447
$ws.Cells.Item(19, 14).Value = 1.012997767402241

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.020611839641401
$ws.Cells.Item(20, 4).Value = 1.023532627453575
$ws.Cells.Item(20, 5).Value = 1.021619282286715
$ws.Cells.Item(20, 6).Value = 1.017088738400261
$ws.Cells.Item(20, 9).Value = 1.028689223167266
$ws.Cells.Item(20, 10).Value = 1.026913381783489
$ws.Cells.Item(20, 11).Value = 1.026950307808486
$ws.Cells.Item(20, 12).Value = 1.025043841935285
$ws.Cells.Item(20, 13).Value = 1.020529701277858
$ws.Cells.Item(20, 14).Value = 1.012877160845895

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.019051839521417
$ws.Cells.Item(21, 4).Value = 1.022074401208775
$ws.Cells.Item(21, 5).Value = 1.020283601184139
$ws.Cells.Item(21, 6).Value = 1.014795108975657
$ws.Cells.Item(21, 9).Value = 1.028422749804144
$ws.Cells.Item(21, 10).Value = 1.025762052469119
$ws.Cells.Item(21, 11).Value = 1.025707147673053
$ws.Cells.Item(21, 12).Value = 1.023923228548916
$ws.Cells.Item(21, 13).Value = 1.018455993070097
$ws.Cells.Item(21, 14).Value = 1.01248363057664

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.018068796964824
$ws.Cells.Item(22, 4).Value = 1.02115574663037
$ws.Cells.Item(22, 5).Value = 1.019442254338921
$ws.Cells.Item(22, 6).Value = 1.013349985543414
$ws.Cells.Item(22, 9).Value = 1.0282528728755
$ws.Cells.Item(22, 10).Value = 1.025035611416291
$ws.Cells.Item(22, 11).Value = 1.024923162581056
$ws.Cells.Item(22, 12).Value = 1.023216521209079
$ws.Cells.Item(22, 13).Value = 1.017148819648123
$ws.Cells.Item(22, 14).Value = 1.01223509702887

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.018590138913081
$ws.Cells.Item(23, 4).Value = 1.021642917061266
$ws.Cells.Item(23, 5).Value = 1.01988841758673
$ws.Cells.Item(23, 6).Value = 1.014116366647784
$ws.Cells.Item(23, 9).Value = 1.028343149808006
$ws.Cells.Item(23, 10).Value = 1.025420956110326
$ws.Cells.Item(23, 11).Value = 1.025338994199423
$ws.Cells.Item(23, 12).Value = 1.023591364534603
$ws.Cells.Item(23, 13).Value = 1.017842100806662
$ws.Cells.Item(23, 14).Value = 1.012366954921438

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.020637956510274
$ws.Cells.Item(24, 4).Value = 1.023557044730983
$ws.Cells.Item(24, 5).Value = 1.021641649383368
$ws.Cells.Item(24, 6).Value = 1.017127141601016
$ws.Cells.Item(24, 9).Value = 1.028693651202271
$ws.Cells.Item(24, 10).Value = 1.026932641200102
$ws.Cells.Item(24, 11).Value = 1.026971110084944
$ws.Cells.Item(24, 12).Value = 1.025062593491828
$ws.Cells.Item(24, 13).Value = 1.020564411842762
$ws.Cells.Item(24, 14).Value = 1.012883739865426

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.0230050310077
$ws.Cells.Item(25, 4).Value = 1.025770671691056
$ws.Cells.Item(25, 5).Value = 1.023669655818002
$ws.Cells.Item(25, 6).Value = 1.020608477350265
$ws.Cells.Item(25, 9).Value = 1.02909028971857
$ws.Cells.Item(25, 10).Value = 1.028675982699109
$ws.Cells.Item(25, 11).Value = 1.028855057520654
$ws.Cells.Item(25, 12).Value = 1.026760801756369
$ws.Cells.Item(25, 13).Value = 1.023709527329995
$ws.Cells.Item(25, 14).Value = 1.013478703371661
